$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 13610.167
$ws.Range("I21").Value = 8348
$ws.Range("K21").Value = 8348
$ws.Range("M21").Value = -7880

$ws.Range("H23").Value = 13610.167
$ws.Range("I23").Value = 8348
$ws.Range("K23").Value = 8348
$ws.Range("M23").Value = -8114

$ws.Range("H137").Value = 1349.9
$ws.Range("I137").Value = 999.6667
$ws.Range("J137").Value = 1500
$ws.Range("K137").Value = 2999.0001
$ws.Range("L137").Value = 4500
$ws.Range("M137").Value = -449.0001000000002
$ws.Range("N137").Value = -9600

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1833.25
$ws.Range("I61").Value = 1833.25
$ws.Range("K61").Value = 1833.25
$ws.Range("M61").Value = -1621.25

$ws.Range("H63").Value = 10453
$ws.Range("I63").Value = 900
$ws.Range("J63").Value = 20006
$ws.Range("K63").Value = 900
$ws.Range("L63").Value = 20006
$ws.Range("M63").Value = -214
$ws.Range("N63").Value = -21378

$ws.Range("H66").Value = 10453
$ws.Range("I66").Value = 900
$ws.Range("J66").Value = 20006
$ws.Range("K66").Value = 4500
$ws.Range("L66").Value = 100030
$ws.Range("M66").Value = -1068
$ws.Range("N66").Value = -106894

$ws.Range("H74").Value = 8866.615
$ws.Range("I74").Value = 8772.25
$ws.Range("K74").Value = 8772.25
$ws.Range("M74").Value = -7898.25

$ws.Range("H77").Value = 8866.615
$ws.Range("I77").Value = 8772.25
$ws.Range("K77").Value = 43861.25
$ws.Range("M77").Value = -39493.25

$ws.Range("H136").Value = 1833.25
$ws.Range("I136").Value = 1833.25
$ws.Range("K136").Value = 5499.75
$ws.Range("M136").Value = -2949.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 964.6
$ws.Range("I80").Value = 655.6667
$ws.Range("K80").Value = 655.6667
$ws.Range("M80").Value = 342.3333

$ws.Range("H83").Value = 964.6
$ws.Range("I83").Value = 655.6667
$ws.Range("K83").Value = 3278.3335
$ws.Range("M83").Value = 1713.6665

$ws.Range("H94").Value = 3825.125
$ws.Range("I94").Value = 2720.4
$ws.Range("J94").Value = 5666.3335
$ws.Range("K94").Value = 2720.4
$ws.Range("L94").Value = 5666.3335
$ws.Range("M94").Value = -2269.4
$ws.Range("N94").Value = -6568.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2474.75
$ws.Range("I58").Value = 2733
$ws.Range("K58").Value = 2733
$ws.Range("M58").Value = -2530

$ws.Range("H136").Value = 2474.75
$ws.Range("I136").Value = 2733
$ws.Range("K136").Value = 8199
$ws.Range("M136").Value = -5649

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

$ws.Range("H80").Value = 9597.5
$ws.Range("I80").Value = 2845
$ws.Range("K80").Value = 2845
$ws.Range("M80").Value = -1847

$ws.Range("H83").Value = 9597.5
$ws.Range("I83").Value = 2845
$ws.Range("K83").Value = 14225
$ws.Range("M83").Value = -9233

$ws.Range("H97").Value = 9999.5
$ws.Range("I97").Value = 9999.5
$ws.Range("K97").Value = 9999.5
$ws.Range("M97").Value = -9503.5

$ws.Range("H102").Value = 6320.5
$ws.Range("I102").Value = 5464.6665
$ws.Range("K102").Value = 5464.6665
$ws.Range("M102").Value = -3842.6665

$ws.Range("H122").Value = 3951
$ws.Range("I122").Value = 4007.375
$ws.Range("K122").Value = 12022.125
$ws.Range("M122").Value = -9572.125

$ws.Range("H132").Value = 4583
$ws.Range("I132").Value = 4500
$ws.Range("J132").Value = 4666
$ws.Range("K132").Value = 13500
$ws.Range("L132").Value = 13998
$ws.Range("M132").Value = -10970
$ws.Range("N132").Value = -19058

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 6500
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 6500
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 6500
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -6726

$ws.Range("H28").Value = 6500
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 6500
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 6500
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value = -6964

$ws.Range("H37").Value = 6500
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 6500
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 6500
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -6714

$ws.Range("H55").Value = 1805.9375
$ws.Range("I55").Value = 4144.6
$ws.Range("J55").Value = 742.9091
$ws.Range("K55").Value = 4144.6
$ws.Range("L55").Value = 742.9091
$ws.Range("M55").Value = -3971.6
$ws.Range("N55").Value = -1088.9091

$ws.Range("H122").Value = 4555.7144
$ws.Range("I122").Value = 4555.7144
$ws.Range("K122").Value = 13667.1432
$ws.Range("M122").Value = -11217.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 18812.5
$ws.Range("I2").Value = 19083.334
$ws.Range("K2").Value = 19083.334
$ws.Range("M2").Value = -18971.334

$ws.Range("H4").Value = 18888.334
$ws.Range("I4").Value = 18666
$ws.Range("J4").Value = 19333
$ws.Range("K4").Value = 18666
$ws.Range("L4").Value = 19333
$ws.Range("M4").Value = -18553
$ws.Range("N4").Value = -19559

$ws.Range("H81").Value = 2015.875
$ws.Range("I81").Value = 1875
$ws.Range("J81").Value = 3002
$ws.Range("K81").Value = 3750
$ws.Range("L81").Value = 6004
$ws.Range("M81").Value = -2689
$ws.Range("N81").Value = -8126

$ws.Range("H84").Value = 2015.875
$ws.Range("I84").Value = 1875
$ws.Range("J84").Value = 3002
$ws.Range("K84").Value = 18750
$ws.Range("L84").Value = 30020
$ws.Range("M84").Value = -13446
$ws.Range("N84").Value = -40628

$ws.Range("H96").Value = 3559
$ws.Range("J96").Value = 1114
$ws.Range("L96").Value = 1114
$ws.Range("N96").Value = -3860

$ws.Range("H126").Value = 1317.3214
$ws.Range("I126").Value = 1093.125
$ws.Range("K126").Value = 3279.375
$ws.Range("M126").Value = -809.375

$ws.Range("H136").Value = 4704
$ws.Range("I136").Value = 4704
$ws.Range("K136").Value = 14112
$ws.Range("M136").Value = -11562
